# ISS-40: Add processing solid fill
#
# Adds a new red-filled rectangle ("Прямоугольник 1") to slide 2, placed
# after the two existing p:graphicFrame OLE objects in the shape tree.
#
# The new shape needs PowerPoint's normal "quick style" block
# (p:style -> lnRef/fillRef/effectRef/fontRef, accent1 scheme) that the
# UI attaches automatically when you draw a shape. AddShape() in this
# runtime does not synthesize that block, so instead we copy an existing
# standalone shape that already carries the accent1 style (slide 3's
# "Прямоугольник 3") and paste it onto slide 2, then restyle it in place
# (name / position / size / solid red fill) to match the target shape.

$p = $ppt.ActivePresentation

$targetSlide = $p.Slides.Item(2)
$styleSourceSlide = $p.Slides.Item(3)

# Donor shape: standalone (not grouped) rectangle that already has the
# accent1-based p:style + centered, run-less txBody we need.
$styleDonor = $styleSourceSlide.Shapes.Item("Прямоугольник 3")
$styleDonor.Copy()

$pasted = $targetSlide.Shapes.Paste()
$newShape = $pasted.Item(1)

$newShape.Name = "Прямоугольник 1"

# Shape.Left/Top/Width/Height are in points; source positions are EMU.
$emuPerPoint = 12700
$newShape.Left = 6007894 / $emuPerPoint
$newShape.Top = 2871788 / $emuPerPoint
$newShape.Width = 1064419 / $emuPerPoint
$newShape.Height = 785812 / $emuPerPoint

# Solid red fill (also clears the donor's picture fill).
$newShape.Fill.Solid()
$newShape.Fill.ForeColor.RGB = 255

Write-Output "Added shape '$($newShape.Name)' to slide 2 (shapes now: $($targetSlide.Shapes.Count))"
